# Daily-push data update: a new sampling row for 2026/01/30 (Fri) was
# recorded and needs to be inserted in date order, right before the
# existing 2026/12/29 block (row 737), pushing rows 737:778 down to
# 738:779 and extending the sheet's used range to A1:D779.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 737 down by one row.
$ws.Rows("737").Insert()

# Column A holds plain date-as-text ("2026/01/30"), not a real date
# serial. A leading apostrophe forces Excel to keep it as literal text
# instead of auto-converting to a date value; ClearFormats() then drops
# the transient quote-prefix/number-format styling so the cell ends up
# with no explicit style, matching the rest of the column.
$ws.Cells.Item(737, 1).Value = "'2026/01/30"
$ws.Cells.Item(737, 1).ClearFormats()
$ws.Cells.Item(737, 2).Value = "金"
$ws.Cells.Item(737, 3).Value = 2
$ws.Cells.Item(737, 4).Value = 201
